# Populate the "খাতা/পত্রের সংখ্যা" (quantity) column G for the three
# billing line items that were left blank, so their dependent amount
# formulas in column I (and the grand-total SUM in I32) recalculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
